$d = $word.ActiveDocument

# The document contains five occurrences of an <id>...</id> tag, each
# currently split across three separate runs:
#   run1: "<id>"        (Courier New, color 7f6000)
#   run2: "p154r_N"      (color 000000, default font)
#   run3: "</id>"       (Courier New, color 7f6000)
#
# The edit merges each triple into a single run "<id>p154r_N</id>" that
# keeps the Courier New / 7f6000 formatting of the surrounding tag runs.
# Using Find & Replace across the whole (already-contiguous) text achieves
# exactly this merge, taking on the formatting of the run at the start of
# the matched range.

for ($i = 1; $i -le 5; $i++) {
    $old = "<id>p154r_$i</id>"
    $new = "<id>p154r_$i</id>"
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}
